$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11; this shifts existing rows 11-13 down to 12-14
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with this week's data.
$ws.Range("A11").Value = 12
$ws.Range("B11").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C11").Value = "Metropolitana"
$ws.Range("D11").Value = 44460
$ws.Range("E11").Value = 13
$ws.Range("F11").Value = 100112021
$ws.Range("G11").Value = "Ají"
$ws.Range("H11").Value = "Americana (o)"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 30
$ws.Range("K11").Value = 95000
$ws.Range("L11").Value = 95000
$ws.Range("M11").Value = 95000
$ws.Range("N11").Value = "$/caja 25 kilos"
$ws.Range("O11").Value = "Provincia de Limarí"
$ws.Range("P11").Value = 3800
$ws.Range("Q11").Value = 25
$ws.Range("R11").Value = "Hortaliza"
